$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended to the list (rows 27-34)
$newRows = @(
    @('ELIABET VARELA REYES', '1212126556', 'eliabetvarela@gmail.com', '13589', 'CONSULTA', '$500', '17:52', '001-20220906V', '06/09/2022'),
    @('PRUEBA LANDSCAPE', '1238923154', 'aaa@.com', '13223', 'SERVICIO', '$132', '18:00', '002-20220906V', '06/09/2022'),
    @('TEST LANDSCAPE', '12312321', 'aaa@.com', '12378', 'CONSULTA', '$123', '18:07', '003-20220906V', '06/09/2022'),
    @('TESTDOCUMENTLANDSCAPE', '1', 'a', '12332', 'SERVICIO|', '$1', '18:09', '004-20220906V', '06/09/2022'),
    @('CARLOS PEÑA', '213712373127', 'cpena@gmail.com', '31265', 'CONSULTA', '$300', '18:12', '005-20220906V', '06/09/2022'),
    @('NOMBRE', '213132123', 'correo', '12323', 'SERVICIO', '$500', '04:41', '001-20220906N', '06/09/2022'),
    @('NOMBRE', '2132345454', 'correo', '21323', 'SERVICIO', '$500', '04:43', '001-20220906M', '06/09/2022'),
    @('PRUEBA CARPETAS', '2183981233', 'asdasd', '98123', 'SERVICIO', '$300', '04:46', '001-20220906V', '06/09/2022')
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowValues[$c - 1]
        $cell.Style = "Normal"
    }
}
